# Arbeitszeit Pichler - add two new work-log entries (GitHub Repo erstellt /
# Code-Basis erstellt) and widen column D to fit the new dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: 01.02.2019 / 1 Stunden / GitHub Repo erstellt -------------------
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats - reuse the date style
$ws.Range("D8").Value = 43497
$ws.Range("E8").Value = "1 Stunden"
$ws.Range("F8").Value = "GitHub Repo erstellt"

# --- Row 9: 04.02.2019 / 1 Stunden / Code-Basis erstellt ---------------------
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats - reuse the date style
$ws.Range("D9").Value = 43500
$ws.Range("E9").Value = "1 Stunden"
$ws.Range("F9").Value = "Code-Basis erstellt"

# --- Column D needs to be wide enough for the "Tätigkeit" dates/labels ------
$ws.Columns("D").ColumnWidth = 21

# --- Move the active selection to D10, like in the saved workbook ----------
$ws.Range("D10").Select() | Out-Null
